{"js": "// Apply the \"carta_generica\" letter edits: switch the student references\n// from masculine to feminine form, update the day-of-month and the RUT\n// check digit, and change the semester referenced in the letter.\n//\n// Each change is located via a unique, case-sensitive text search on\n// context.document.body so the existing run formatting (bold, etc.) is\n// preserved by replacing only the matched range's text.\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText, options) {\n  const searchOptions = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(searchText, searchOptions);\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${searchText}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Valpara\u00edso, 4 de marzo de 2024\" -> \"Valpara\u00edso, 24 de marzo de 2024\"\n//    (matchWholeWord avoids also matching the \"4\" inside \"2024\")\nawait replaceOnce(\"4\", \"24\", { matchWholeWord: true });\n\n// 2) \"nuestro alumno Se\u00f1or\" -> \"nuestra alumna Se\u00f1orita\"\nawait replaceOnce(\"nuestro alumno Se\u00f1or\", \"nuestra alumna Se\u00f1orita\");\n\n// 3) RUT check digit \"21061253-K\" -> \"21061253-3\"\nawait replaceOnce(\"21061253-K\", \"21061253-3\");\n\n// 4) \", en calidad de alumno\" -> \", en calidad de alumna\"\nawait replaceOnce(\", en calidad de alumno\", \", en calidad de alumna\");\n\n// 5) \"El se\u00f1or\" -> \"La se\u00f1orita\"\nawait replaceOnce(\"El se\u00f1or\", \"La se\u00f1orita\");\n\n// 6) \"Octavo semestre\" -> \"Sexto Semestre\"\nawait replaceOnce(\"Octavo semestre\", \"Sexto Semestre\");\n\n// 7) \"EL ALUMNO\" -> \"LA ALUMNA\"\nawait replaceOnce(\"EL ALUMNO\", \"LA ALUMNA\");\n", "ps1": "# Apply the \"carta_generica\" letter edits: switch the student references\n# from masculine to feminine form, update the day-of-month and the RUT\n# check digit, and change the semester referenced in the letter.\n#\n# Each change is performed with Find/Replace over the whole document\n# range so existing run formatting (bold, etc.) on the matched text is\n# preserved.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text([string]$FindText, [string]$ReplaceText, [bool]$MatchWholeWord) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n\n    $found = $range.Find.Execute(\n        $FindText,\n        $true,\n        $MatchWholeWord,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $ReplaceText,\n        2\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $FindText\"\n    }\n}\n\n# 1) \"Valpara\u00edso, 4 de marzo de 2024\" -> \"Valpara\u00edso, 24 de marzo de 2024\"\n#    (MatchWholeWord avoids also matching the \"4\" inside \"2024\")\nReplace-Text \"4\" \"24\" $true\n\n# 2) \"nuestro alumno Se\u00f1or\" -> \"nuestra alumna Se\u00f1orita\"\nReplace-Text \"nuestro alumno Se\u00f1or\" \"nuestra alumna Se\u00f1orita\" $false\n\n# 3) RUT check digit \"21061253-K\" -> \"21061253-3\"\nReplace-Text \"21061253-K\" \"21061253-3\" $false\n\n# 4) \", en calidad de alumno\" -> \", en calidad de alumna\"\nReplace-Text \", en calidad de alumno\" \", en calidad de alumna\" $false\n\n# 5) \"El se\u00f1or\" -> \"La se\u00f1orita\"\nReplace-Text \"El se\u00f1or\" \"La se\u00f1orita\" $false\n\n# 6) \"Octavo semestre\" -> \"Sexto Semestre\"\nReplace-Text \"Octavo semestre\" \"Sexto Semestre\" $false\n\n# 7) \"EL ALUMNO\" -> \"LA ALUMNA\"\nReplace-Text \"EL ALUMNO\" \"LA ALUMNA\" $false\n"}
